# "Add files via upload / i have updated something inside the file"
#
# The authored change simply adds one more row of data to the bottom of the
# "Pending Task" sheet: a new shared string "I made changes here" is placed
# in cell B14, which pushes the sheet's used range/dimension out to F14 and
# moves the active selection down to B14.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Write the new note into B14 (this naturally grows the shared-string table
# and the sheet dimension to A1:F14).
$ws.Range("B14").Value = "I made changes here"

# Mirror the author's final cursor position (selection moved to B14).
$ws.Range("B14").Select() | Out-Null
